$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.203.94"
$ws.Range("E2").Value = "  -0.56%  "
$ws.Range("D3").Value = "1.585.35"
$ws.Range("E3").Value = "  -0.43%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.54"
$ws.Range("E5").Value = "  +0.72%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.502"
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -0.22%  "
$ws.Range("E9").Value = "  -0.93%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.23"
$ws.Range("E10").Value = "  -2.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0846"
$ws.Range("E11").Value = "  +0.17%  "
$ws.Range("D12").Value = "1.808.47"
$ws.Range("E12").Value = "  -0.37%  "
$ws.Range("D13").Value = "1.592.29"
$ws.Range("E13").Value = "  -0.53%  "
$ws.Range("E14").Value = "  -1.55%  "
$ws.Range("E15").Value = "  -0.41%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.02"
$ws.Range("E16").Value = "  -1.21%  "
$ws.Range("D17").Value = "26.200.76"
$ws.Range("E17").Value = "  -0.52%  "
$ws.Range("E18").Value = "  -0.72%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "213.99"
$ws.Range("E19").Value = "  +0.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.33"
$ws.Range("E20").Value = "  -0.84%  "
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("E22").Value = "  -1.39%  "
$ws.Range("E23").Value = "  -0.79%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.95"
$ws.Range("E24").Value = "  +0.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.72"
$ws.Range("E25").Value = "  -0.53%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("E27").Value = "  -1.09%  "
$ws.Range("E28").Value = "  -1.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.14"
$ws.Range("E29").Value = "  -1.09%  "
$ws.Range("E30").Value = "  -1.78%  "
$ws.Range("E31").Value = "  +0.68%  "
$ws.Range("E32").Value = "  -1.25%  "
$ws.Range("D33").Value = "1.356.70"
$ws.Range("E33").Value = "  +4.76%  "
$ws.Range("E34").Value = "  -1.76%  "
$ws.Range("E35").Value = "  -0.29%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.46"
$ws.Range("E36").Value = "  -1.29%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.580"
$ws.Range("E37").Value = "  -5.18%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0167"
$ws.Range("E38").Value = "  -0.63%  "
$ws.Range("E39").Value = "  +0.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.81"
$ws.Range("E40").Value = "  +2.98%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  -0.09%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.768"
$ws.Range("E42").Value = "  +0.76%  "
$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.926"
$ws.Range("E43").Value = "  -17.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.14"
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("D45").Value = "1.721.14"
$ws.Range("E45").Value = "  -0.34%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "60.81"
$ws.Range("E46").Value = "  -3.12%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "86.12"
$ws.Range("E47").Value = "  -2.91%  "
$ws.Range("E48").Value = "  -0.53%  "
$ws.Range("E49").Value = "  -2.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0979"
$ws.Range("E50").Value = "  -2.10%  "
$ws.Range("E51").Value = "  -1.10%  "
